$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The checklist currently ends at row 36 (year 1987). Insert 24 fresh rows
# right after it (12 months x 2 years: 1988 and 1989), each one inheriting
# the formatting of the row above it (matches how Excel's own "insert row"
# carries formatting down).
for ($i = 0; $i -lt 24; $i++) {
    $ws.Rows.Item(37).Insert()
}

$row = 37
foreach ($year in 1988, 1989) {
    for ($month = 1; $month -le 12; $month++) {
        $ws.Range("A$row").Value = $year
        $ws.Range("B$row").Value = "{0}/{1:D2}" -f $year, $month
        $ws.Range("D$row").Value = "Kadokawa Shoten"
        $row++
    }
}

# Fill the image-name formula down across all the new rows in one shot so it
# forms a single shared-formula group, same as the existing E4:E36 block.
$ws.Range("E37:E60").Formula = "=CONCATENATE(SUBSTITUTE(B37,""/"",""""), "".jpg"")"

# Leave the 1989 block selected, matching where editing finished.
$ws.Range("E49:E60").Select()
